# Update match data for "Paraguay Division Profesional" workbook.
#
# The underlying data rows (identified by the sequential "id" column A,
# which must stay untouched) were re-ordered / corrected for a handful of
# fixtures that share the same kickoff date. We reproduce this by reading
# the affected rows (columns B..AD, i.e. everything except the id column)
# as value arrays and writing them back in the corrected order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($sheet, $row1, $row2) {
    $range1 = $sheet.Range("B$row1`:AD$row1")
    $range2 = $sheet.Range("B$row2`:AD$row2")
    $v1 = $range1.Value2
    $v2 = $range2.Value2
    $range1.Value2 = $v2
    $range2.Value2 = $v1
}

# Rows 2, 3, 4 (all dated 2023-06-07 85416...) get cyclically rotated:
#   new row 2 <- old row 4
#   new row 3 <- old row 2
#   new row 4 <- old row 3
$r2 = $ws.Range("B2:AD2").Value2
$r3 = $ws.Range("B3:AD3").Value2
$r4 = $ws.Range("B4:AD4").Value2

$ws.Range("B2:AD2").Value2 = $r4
$ws.Range("B3:AD3").Value2 = $r2
$ws.Range("B4:AD4").Value2 = $r3

# Rows 134 and 135 (both dated 2023-11-10) swap places entirely.
Swap-Rows $ws 134 135

# Rows 143 and 145 (both dated 2023-11-18) swap places entirely; row 144
# in between is left untouched.
Swap-Rows $ws 143 145
